$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = $cell.Value2
        if ($current -eq 0) {
            $cell.Value2 = 1
        } else {
            $cell.Value2 = 0
        }
    }
}

$ws.Range("U8").Select()
